$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B6").Value = 'Sie haben {{num_correct}} aus {{num_items}} Namen richtig erkannt ({{perc_correct}}%,\\ {{FP}} falsch gewählt\\ und {{FN}} nicht erkannt.\\Das ergibt **{{points}}/100** Punkte.'
$ws.Range("C6").Value = 'You recognized {{num_correct}} out of {{num_items}} names correctly ({{perc_correct}}%),\\ you assigned {{FP}} wrongly,\\ and you missed {{FN}}.\\ This  yields **{{points}}/100** points.'

$ws.Range("B6").Select()
